# Remove the stray duplicated header row (row 95) from "DEVAM EDEN";
# all rows below it (96..139) shift up by one (95..138).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("DEVAM EDEN")
$ws1.Rows.Item(95).Delete()

# Restore the view/selection state recorded in the saved workbook.
$ws2 = $wb.Worksheets.Item("İHALEDE")
$ws2.Activate()
[void]$ws2.Range("F54").Select()

$ws1.Activate()
[void]$ws1.Range("I84").Select()
